# Regenerate save_data to use K (strikeouts) instead of Strike# (pitch-count based
# strike total), recalculated std/mean, and write the newly computed s_vals (column G)
# back into the worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G ("K") values for rows 2..32, replacing the old Strike# values.
$kValues = @(5, 6, 12, 1, 4, 4, 6, 5, 5, 4, 7, 0, 8, 3, 4, 9, 8, 7, 3, 3, 7, 4, 7, 6, 11, 3, 8, 6, 4, 3, 3)

$row = 2
foreach ($val in $kValues) {
    $ws.Cells.Item($row, 7).Value = $val
    $row++
}
